$wb = $excel.ActiveWorkbook

function Set-Rows {
    param(
        $ws,
        [int]$startRow,
        $rows
    )
    $r = $startRow
    foreach ($rowData in $rows) {
        $c = 1
        foreach ($val in $rowData) {
            $ws.Cells.Item($r, $c).Value = $val
            $c++
        }
        $r++
    }
}

# ---------------------------------------------------------------------------
# Sheet "Home win": swap rows 2/3 with new match data and append rows 6-12
# ---------------------------------------------------------------------------
$wsHome = $wb.Worksheets.Item("Home win")

Set-Rows $wsHome 2 @(
    , @("29-01-2025 22:00", "BRAZIL", "GAÚCHO - 1", "Ypiranga-RS - São Luiz", 86.7, 2.25)
    , @("29-01-2025 21:00", "EL-SALVADOR", "PRIMERA DIVISION", "Cacahuatique - Fuerte San Francisco", 73.3, 1.77)
)

Set-Rows $wsHome 6 @(
    , @("30-01-2025 20:00", "WORLD", "UEFA EUROPA LEAGUE", "Dynamo Kyiv - Rīgas FS", 73.3, 1.85)
    , @("30-01-2025 20:00", "WORLD", "UEFA EUROPA LEAGUE", "Maccabi Tel Aviv - FC Porto", 70, 7)
    , @("30-01-2025 20:00", "WORLD", "UEFA EUROPA LEAGUE", "Rangers - Union St. Gilloise", 73.3, 2)
    , @("30-01-2025 20:00", "WORLD", "UEFA EUROPA LEAGUE", "AS Roma - Eintracht Frankfurt", 80, 1.73)
    , @("30-01-2025 23:00", "BRAZIL", "GAÚCHO - 1", "Pelotas - Avenida", 73.3, 2.1)
    , @("30-01-2025 00:30", "BRAZIL", "GAÚCHO - 1", "Caxias - Brasil DE Pelotas", 71.7, 1.7)
    , @("30-01-2025 08:30", "INDONESIA", "LIGA 2", "Persikota Tangerang - Sriwijaya FC", 80, 1.91)
)

# ---------------------------------------------------------------------------
# Sheet "Btts": replace row 11 and append rows 12-13
# ---------------------------------------------------------------------------
$wsBtts = $wb.Worksheets.Item("Btts")

Set-Rows $wsBtts 11 @(
    , @("30-01-2025 20:00", "WORLD", "UEFA EUROPA LEAGUE", "FCSB - Manchester United", 76, 1.8)
    , @("30-01-2025 21:30", "BRAZIL", "CARIOCA - 1", "Boavista SC - Madureira", 75.8, 1.91)
    , @("30-01-2025 01:00", "MEXICO", "LIGA MX", "Club Queretaro - Pachuca", 76.7, 1.7)
)

# ---------------------------------------------------------------------------
# Sheet "Over_Under": append rows 8-9
# ---------------------------------------------------------------------------
$wsOU = $wb.Worksheets.Item("Over_Under")

Set-Rows $wsOU 8 @(
    , @("30-01-2025 20:00", "WORLD", "UEFA EUROPA LEAGUE", "AS Roma - Eintracht Frankfurt", 80, 1.8, 53.3, 3)
    , @("30-01-2025 20:00", "WORLD", "UEFA EUROPA LEAGUE", "SC Braga - Lazio", 86.7, 1.73, 33.3, 2.75)
)

Write-Output "updated predictions sheets"
